# Weekly update: insert two new price records (2021-10-05) at the top of the
# Femacal de La Calera - Alcachofa data block, shifting the existing history
# down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 172; this pushes the former rows 172:198
# down to 174:200 and keeps their contents intact (including styles/number
# formats), matching the "older week pushed down" pattern seen in the diff.
$ws.Rows("172:173").Insert()

# New row 172: Alcachofa, Argentina(o), Primera
$ws.Range("A172").Value = 3
$ws.Range("B172").Value = "Femacal de La Calera"
$ws.Range("C172").Value = "Coquimbo"
$ws.Range("D172").Value = 44474
$ws.Range("E172").Value = 5
$ws.Range("F172").Value = 100112013
$ws.Range("G172").Value = "Alcachofa"
$ws.Range("H172").Value = "Argentina(o)"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 130
$ws.Range("K172").Value = 9000
$ws.Range("L172").Value = 9500
$ws.Range("M172").Value = 9269
$ws.Range("N172").Value = "$/caja 50 unidades"
$ws.Range("O172").Value = "Provincia de Limarí"
$ws.Range("P172").Value = 185
$ws.Range("Q172").Value = 50
$ws.Range("R172").Value = "Hortaliza"

# New row 173: Alcachofa, Española, Extra
$ws.Range("A173").Value = 3
$ws.Range("B173").Value = "Femacal de La Calera"
$ws.Range("C173").Value = "Coquimbo"
$ws.Range("D173").Value = 44474
$ws.Range("E173").Value = 5
$ws.Range("F173").Value = 100112013
$ws.Range("G173").Value = "Alcachofa"
$ws.Range("H173").Value = "Española"
$ws.Range("I173").Value = "Extra"
$ws.Range("J173").Value = 130
$ws.Range("K173").Value = 9000
$ws.Range("L173").Value = 9500
$ws.Range("M173").Value = 9269
$ws.Range("N173").Value = "$/caja 30 unidades"
$ws.Range("O173").Value = "Provincia de Limarí"
$ws.Range("P173").Value = 309
$ws.Range("Q173").Value = 30
$ws.Range("R173").Value = "Hortaliza"
